$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Bg_Jimroom.jpg now has a source and is a Placeholder instead of Not Done ---
$ws.Range("E8").Value = "Image by Kyle"
$ws.Range("F8").Value = "Placeholder"

# --- Insert 9 new rows for the intro-sequence item placeholders (new rows 18-26) ---
$ws.Rows("18:26").Insert()

$ws.Range("A18").Value = "image"
$ws.Range("B18").Value = "Item_Phone.png"
$ws.Range("C18").Value = "Phone item, used in Jim's Room"
$ws.Range("D18").Value = "460 X 856"
$ws.Range("E18").Value = "Image by Kyle"
$ws.Range("F18").Value = "Placeholder"

$ws.Range("A19").Value = "image"
$ws.Range("B19").Value = "Item_Clock.png"
$ws.Range("C19").Value = "Clock, shows player how much time they have"
$ws.Range("F19").Value = "Not Done"

$ws.Range("A20").Value = "image"
$ws.Range("B20").Value = "Item_Splicer.png"
$ws.Range("C20").Value = "DNA Splicer, Madam Feline's item"
$ws.Range("F20").Value = "Not Done"

$ws.Range("A21").Value = "image"
$ws.Range("B21").Value = "Item_Flour.png"
$ws.Range("C21").Value = "Bag of flour, Chris's item"
$ws.Range("F21").Value = "Not Done"

$ws.Range("A22").Value = "image"
$ws.Range("B22").Value = "Item_Fleece.png"
$ws.Range("C22").Value = "Golden fleece blanket, Jason's item"
$ws.Range("F22").Value = "Not Done"

$ws.Range("A23").Value = "image"
$ws.Range("B23").Value = "Item_Football.png"
$ws.Range("C23").Value = "Deflated football, Coach Dave's item"
$ws.Range("F23").Value = "Not Done"

$ws.Range("A24").Value = "image"
$ws.Range("B24").Value = "Item_Watch.png"
$ws.Range("C24").Value = "Diamond Pocket Watch, Sir Edmond's item"
$ws.Range("F24").Value = "Not Done"

$ws.Range("A25").Value = "image"
$ws.Range("B25").Value = "Item_Heels.png"
$ws.Range("C25").Value = "High Heels, Kim's item"
$ws.Range("F25").Value = "Not Done"

$ws.Range("A26").Value = "image"
$ws.Range("B26").Value = "Item_Candle.png"
$ws.Range("C26").Value = "Candlestick, Colonel Ketchup's item"
$ws.Range("F26").Value = "Not Done"

# --- New sound placeholder ("Phone ringing") reuses the next blank filler row (row 42) ---
$ws.Range("A42").Value = "sound"
$ws.Range("C42").Value = "Phone ringing"

# --- Update the view state to match where the author was working ---
$ws.Range("E23").Select()
